$d = $word.ActiveDocument

# Map of paragraph index -> expected leading text (sanity check) for each
# answer bullet that needs to be prefixed with "Design: " per the commit.
$targets = @{
  8  = "Mostly"
  9  = "Kymys ei täysin ymmärretävä, tarkoitetaanko tällä sitä miten sähkö on onnistunut?"
  13 = "Nothing special to mention"
  14 = "Omien järjestelmien osalta (sähkönjakelu) vaikuttaa moneen järjestelmään säännöt esim. tarpeet tuplasyötöistä, syötönvaihdot jne."
  18 = "91xx Went well, some minor budget challenges"
  19 = "tarjouspyynnöt ja tarjouskierros monimutkainen ja aikaa vievä prosessi."
  23 = "Mostly ok."
  24 = "PES sähkön osalta kannattaisi tehdä telakan omalla väellä."
  28 = "Internal communication ok. External communication with suppliers mostly ok."
  29 = "Suunnittelua tehdään samaan aikaan monella osastolla"
  33 = "Some things went to correct direction but regarding TK the opposite way."
  34 = "protolaivat olivat vaikeita tehdä"
}

# Insert before paragraphs from the bottom up so earlier insertions do not
# shift the indices of paragraphs that still need to be processed.
$indices = $targets.Keys | Sort-Object -Descending

foreach ($idx in $indices) {
    $para = $d.Paragraphs.Item($idx)
    $rng = $para.Range
    if ($rng.Text.StartsWith($targets[$idx])) {
        $rng.InsertBefore("Design: ")
    } else {
        Write-Host ("WARNING: paragraph " + $idx + " did not match expected text")
    }
}

$d.Save()
